# Regenerate merged AHB files
# - rename the "_old"/"_new" suffixed header labels to "_FV2410"/"_FV2504"
# - wrap the data range in an Excel Table (ListObject)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $labels.Count; $i++) {
    # columns A..J -> "_old" becomes "_FV2410"
    $ws.Cells.Item(1, $i + 1).Value = "$($labels[$i])_FV2410"
    # columns L..U -> "_new" becomes "_FV2504"
    $ws.Cells.Item(1, $i + 12).Value = "$($labels[$i])_FV2504"
}

# Turn the used range into a native Excel table ("Table1")
$rng = $ws.Range("A1:U74")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium9"

# Freeze the header row (split below row 1)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
